# Update NATMI LR-pair stats (Sema4d-Plxnb1) with recomputed TPM-based values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.144900666666667
$ws.Range("H2").Value = 3.434702
$ws.Range("I2").Value = 0.02523133726002265
$ws.Range("J2").Value = 0.02523133726002265
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.480335666666667
$ws.Range("N2").Value = 4.441007
$ws.Range("O2").Value = 0.1826408776454046
$ws.Range("P2").Value = 0.1826408776454046
$ws.Range("Q2").Value = 1.694837291657111
$ws.Range("R2").Value = 15.253535624914
$ws.Range("S2").Value = 0.004608273581337734
$ws.Range("T2").Value = 0.004608273581337734
$ws.Range("G3").Value = 1.144900666666667
$ws.Range("H3").Value = 3.434702
$ws.Range("I3").Value = 0.02523133726002265
$ws.Range("J3").Value = 0.02523133726002265
$ws.Range("O3").Value = 0.07762443032771463
$ws.Range("P3").Value = 0.07762443032771463
$ws.Range("Q3").Value = 0.7203249401728887
$ws.Range("R3").Value = 6.482924461555999
$ws.Range("S3").Value = 0.001958568181215698
$ws.Range("T3").Value = 0.001958568181215698
$ws.Range("G4").Value = 1.144900666666667
$ws.Range("H4").Value = 3.434702
$ws.Range("I4").Value = 0.02523133726002265
$ws.Range("J4").Value = 0.02523133726002265
$ws.Range("M4").Value = 5.953764333333333
$ws.Range("N4").Value = 17.861293
$ws.Range("O4").Value = 0.734563631492074
$ws.Range("P4").Value = 0.734563631492074
$ws.Range("Q4").Value = 6.816468754409555
$ws.Range("R4").Value = 61.34821878968599
$ws.Range("S4").Value = 0.01853402272512351
$ws.Range("T4").Value = 0.01853402272512351
$ws.Range("G5").Value = 1.144900666666667
$ws.Range("H5").Value = 3.434702
$ws.Range("I5").Value = 0.02523133726002265
$ws.Range("J5").Value = 0.02523133726002265
$ws.Range("M5").Value = 0.04191233333333333
$ws.Range("N5").Value = 0.125737
$ws.Range("O5").Value = 0.005171060534806686
$ws.Range("P5").Value = 0.005171060534806686
$ws.Range("Q5").Value = 0.04798545837488888
$ws.Range("R5").Value = 0.4318691253739999
$ws.Range("S5").Value = 0.0001304727723457006
$ws.Range("T5").Value = 0.0001304727723457006
$ws.Range("I6").Value = 0.03646539869776051
$ws.Range("J6").Value = 0.03646539869776051
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.480335666666667
$ws.Range("N6").Value = 4.441007
$ws.Range("O6").Value = 0.1826408776454046
$ws.Range("P6").Value = 0.1826408776454046
$ws.Range("Q6").Value = 2.449450733871
$ws.Range("R6").Value = 22.045056604839
$ws.Range("S6").Value = 0.006660072421848574
$ws.Range("T6").Value = 0.006660072421848573
$ws.Range("I7").Value = 0.03646539869776051
$ws.Range("J7").Value = 0.03646539869776051
$ws.Range("O7").Value = 0.07762443032771463
$ws.Range("P7").Value = 0.07762443032771463
$ws.Range("S7").Value = 0.002830605800586646
$ws.Range("T7").Value = 0.002830605800586646
$ws.Range("I8").Value = 0.03646539869776051
$ws.Range("J8").Value = 0.03646539869776051
$ws.Range("M8").Value = 5.953764333333333
$ws.Range("N8").Value = 17.861293
$ws.Range("O8").Value = 0.734563631492074
$ws.Range("P8").Value = 0.734563631492074
$ws.Range("Q8").Value = 9.851449738028998
$ws.Range("R8").Value = 88.66304764226099
$ws.Range("S8").Value = 0.02678615569123331
$ws.Range("T8").Value = 0.02678615569123331
$ws.Range("I9").Value = 0.03646539869776051
$ws.Range("J9").Value = 0.03646539869776051
$ws.Range("M9").Value = 0.04191233333333333
$ws.Range("N9").Value = 0.125737
$ws.Range("O9").Value = 0.005171060534806686
$ws.Range("P9").Value = 0.005171060534806686
$ws.Range("Q9").Value = 0.06935061956099998
$ws.Range("R9").Value = 0.624155576049
$ws.Range("S9").Value = 0.0001885647840919805
$ws.Range("T9").Value = 0.0001885647840919805
$ws.Range("G10").Value = 3.191626333333333
$ws.Range("H10").Value = 9.574878999999999
$ws.Range("I10").Value = 0.07033710676294723
$ws.Range("J10").Value = 0.07033710676294723
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.480335666666667
$ws.Range("N10").Value = 4.441007
$ws.Range("O10").Value = 0.1826408776454046
$ws.Range("P10").Value = 0.1826408776454046
$ws.Range("Q10").Value = 4.724678295905889
$ws.Range("R10").Value = 42.52210466315299
$ws.Range("S10").Value = 0.0128464309102232
$ws.Range("T10").Value = 0.0128464309102232
$ws.Range("G11").Value = 3.191626333333333
$ws.Range("H11").Value = 9.574878999999999
$ws.Range("I11").Value = 0.07033710676294723
$ws.Range("J11").Value = 0.07033710676294723
$ws.Range("O11").Value = 0.07762443032771463
$ws.Range("P11").Value = 0.07762443032771463
$ws.Range("Q11").Value = 2.008041496129111
$ws.Range("R11").Value = 18.072373465162
$ws.Range("S11").Value = 0.005459877843373422
$ws.Range("T11").Value = 0.005459877843373422
$ws.Range("G12").Value = 3.191626333333333
$ws.Range("H12").Value = 9.574878999999999
$ws.Range("I12").Value = 0.07033710676294723
$ws.Range("J12").Value = 0.07033710676294723
$ws.Range("M12").Value = 5.953764333333333
$ws.Range("N12").Value = 17.861293
$ws.Range("O12").Value = 0.734563631492074
$ws.Range("P12").Value = 0.734563631492074
$ws.Range("Q12").Value = 19.00219102872744
$ws.Range("R12").Value = 171.019719258547
$ws.Range("S12").Value = 0.05166708057243623
$ws.Range("T12").Value = 0.05166708057243623
$ws.Range("G13").Value = 3.191626333333333
$ws.Range("H13").Value = 9.574878999999999
$ws.Range("I13").Value = 0.07033710676294723
$ws.Range("J13").Value = 0.07033710676294723
$ws.Range("M13").Value = 0.04191233333333333
$ws.Range("N13").Value = 0.125737
$ws.Range("O13").Value = 0.005171060534806686
$ws.Range("P13").Value = 0.005171060534806686
$ws.Range("Q13").Value = 0.1337685067581111
$ws.Range("R13").Value = 1.203916560823
$ws.Range("S13").Value = 0.0003637174369143609
$ws.Range("T13").Value = 0.0003637174369143609
$ws.Range("G14").Value = 39.384953
$ws.Range("H14").Value = 118.154859
$ws.Range("I14").Value = 0.8679661572792696
$ws.Range("J14").Value = 0.8679661572792696
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 1.480335666666667
$ws.Range("N14").Value = 4.441007
$ws.Range("O14").Value = 0.1826408776454046
$ws.Range("P14").Value = 0.1826408776454046
$ws.Range("Q14").Value = 58.30295065589033
$ws.Range("R14").Value = 524.7265559030129
$ws.Range("S14").Value = 0.1585261007319951
$ws.Range("T14").Value = 0.1585261007319951
$ws.Range("G15").Value = 39.384953
$ws.Range("H15").Value = 118.154859
$ws.Range("I15").Value = 0.8679661572792696
$ws.Range("J15").Value = 0.8679661572792696
$ws.Range("O15").Value = 0.07762443032771463
$ws.Range("P15").Value = 0.07762443032771463
$ws.Range("Q15").Value = 24.77941077284466
$ws.Range("R15").Value = 223.014696955602
$ws.Range("S15").Value = 0.06737537850253886
$ws.Range("T15").Value = 0.06737537850253886
$ws.Range("G16").Value = 39.384953
$ws.Range("H16").Value = 118.154859
$ws.Range("I16").Value = 0.8679661572792696
$ws.Range("J16").Value = 0.8679661572792696
$ws.Range("M16").Value = 5.953764333333333
$ws.Range("N16").Value = 17.861293
$ws.Range("O16").Value = 0.734563631492074
$ws.Range("P16").Value = 0.734563631492074
$ws.Range("Q16").Value = 234.4887284414096
$ws.Range("R16").Value = 2110.398555972687
$ws.Range("S16").Value = 0.637576372503281
$ws.Range("T16").Value = 0.637576372503281
$ws.Range("G17").Value = 39.384953
$ws.Range("H17").Value = 118.154859
$ws.Range("I17").Value = 0.8679661572792696
$ws.Range("J17").Value = 0.8679661572792696
$ws.Range("M17").Value = 0.04191233333333333
$ws.Range("N17").Value = 0.125737
$ws.Range("O17").Value = 0.005171060534806686
$ws.Range("P17").Value = 0.005171060534806686
$ws.Range("Q17").Value = 1.650715278453666
$ws.Range("R17").Value = 14.856437506083
$ws.Range("S17").Value = 0.004488305541454644
$ws.Range("T17").Value = 0.004488305541454644
